# Apply cryptos list update to sheet1 (D = Price, E = Volume(1h))
# Updated cryptos list on Sat Aug 12 18:45:41 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.435.97"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.96"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.86"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6295"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07676"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.71"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.16"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.031"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6797"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001069"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.184"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.454.55"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.17"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.425"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.82"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1382"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.409"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.352"
$ws.Range("E28").Value = "  +6.30%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05671"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.122"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.030"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.843"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.162"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7084"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.585"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.780"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01789"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.220.22"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.548"
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9123"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.60"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.15"
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000120"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.143"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.006"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.677"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1144"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05716"
$ws.Range("E51").Value = "  +0.11%  "
